$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  34" -> "Volume 30   Number  35"
$ws.Range("A8").Characters(21, 2).Text = "35"

# "Report Covering the Week  8/21/2023  Through  8/27/2023"
# -> "Report Covering the Week  8/28/2023  Through  9/3/2023"
$ws.Range("C9").Characters(27, 9).Text = "8/28/2023"
$ws.Range("C9").Characters(47, 9).Text = "9/3/2023"

# --- Weekly crime-stat table updates (rows 14-30) -------------------------
$ws.Cells.Item(14, 6).Value = 29
$ws.Cells.Item(14, 7).Value = 28
$ws.Cells.Item(14, 8).Value = 3.571428571428
$ws.Cells.Item(14, 9).Value = 274
$ws.Cells.Item(14, 10).Value = 303
$ws.Cells.Item(14, 11).Value = -9.570957095709
$ws.Cells.Item(14, 12).Value = -16.463414634146
$ws.Cells.Item(14, 13).Value = -25.745257452574
$ws.Cells.Item(14, 14).Value = -79.226686884003
$ws.Cells.Item(15, 3).Value = 30
$ws.Cells.Item(15, 4).Value = 37
$ws.Cells.Item(15, 5).Value = -18.918918918918
$ws.Cells.Item(15, 6).Value = 107
$ws.Cells.Item(15, 7).Value = 125
$ws.Cells.Item(15, 8).Value = -14.4
$ws.Cells.Item(15, 9).Value = 979
$ws.Cells.Item(15, 10).Value = 1121
$ws.Cells.Item(15, 11).Value = -12.667261373773
$ws.Cells.Item(15, 12).Value = -3.736479842674
$ws.Cells.Item(15, 13).Value = 9.385474860335
$ws.Cells.Item(15, 14).Value = -56.604609929078
$ws.Cells.Item(16, 3).Value = 319
$ws.Cells.Item(16, 4).Value = 349
$ws.Cells.Item(16, 5).Value = -8.595988538681
$ws.Cells.Item(16, 6).Value = 1400
$ws.Cells.Item(16, 7).Value = 1468
$ws.Cells.Item(16, 8).Value = -4.632152588555
$ws.Cells.Item(16, 9).Value = 11054
$ws.Cells.Item(16, 10).Value = 11730
$ws.Cells.Item(16, 11).Value = -5.763000852514
$ws.Cells.Item(16, 12).Value = 31.877833452636
$ws.Cells.Item(16, 13).Value = -11.617494203246
$ws.Cells.Item(16, 14).Value = -80.488924190274
$ws.Cells.Item(17, 3).Value = 531
$ws.Cells.Item(17, 4).Value = 535
$ws.Cells.Item(17, 5).Value = -0.747663551401
$ws.Cells.Item(17, 6).Value = 2148
$ws.Cells.Item(17, 7).Value = 2134
$ws.Cells.Item(17, 8).Value = 0.656044985941
$ws.Cells.Item(17, 9).Value = 18819
$ws.Cells.Item(17, 10).Value = 17868
$ws.Cells.Item(17, 11).Value = 5.322364002686
$ws.Cells.Item(17, 12).Value = 25.871179185338
$ws.Cells.Item(17, 13).Value = 60.997519034990
$ws.Cells.Item(17, 14).Value = -34.190096516995
$ws.Cells.Item(18, 3).Value = 262
$ws.Cells.Item(18, 4).Value = 321
$ws.Cells.Item(18, 5).Value = -18.380062305296
$ws.Cells.Item(18, 6).Value = 1020
$ws.Cells.Item(18, 7).Value = 1283
$ws.Cells.Item(18, 8).Value = -20.498830865159
$ws.Cells.Item(18, 9).Value = 9370
$ws.Cells.Item(18, 10).Value = 10530
$ws.Cells.Item(18, 11).Value = -11.016144349477
$ws.Cells.Item(18, 12).Value = 19.241537286841
$ws.Cells.Item(18, 13).Value = -23.403907463418
$ws.Cells.Item(18, 14).Value = -86.056755107810
$ws.Cells.Item(19, 3).Value = 962
$ws.Cells.Item(19, 4).Value = 1097
$ws.Cells.Item(19, 5).Value = -12.306289881495
$ws.Cells.Item(19, 6).Value = 4087
$ws.Cells.Item(19, 7).Value = 4286
$ws.Cells.Item(19, 8).Value = -4.643023798413
$ws.Cells.Item(19, 9).Value = 33694
$ws.Cells.Item(19, 10).Value = 34521
$ws.Cells.Item(19, 11).Value = -2.395643231656
$ws.Cells.Item(19, 12).Value = 43.470300191611
$ws.Cells.Item(19, 13).Value = 36.131873459658
$ws.Cells.Item(19, 14).Value = -41.283284539244
$ws.Cells.Item(20, 3).Value = 306
$ws.Cells.Item(20, 4).Value = 284
$ws.Cells.Item(20, 5).Value = 7.746478873239
$ws.Cells.Item(20, 6).Value = 1343
$ws.Cells.Item(20, 7).Value = 1067
$ws.Cells.Item(20, 8).Value = 25.866916588566
$ws.Cells.Item(20, 9).Value = 10632
$ws.Cells.Item(20, 10).Value = 8953
$ws.Cells.Item(20, 11).Value = 18.753490450128
$ws.Cells.Item(20, 12).Value = 64.78611283323
$ws.Cells.Item(20, 13).Value = 53.353526611856
$ws.Cells.Item(20, 14).Value = -85.732114819437
$ws.Cells.Item(21, 3).Value = 2420
$ws.Cells.Item(21, 4).Value = 2633
$ws.Cells.Item(21, 5).Value = -8.089631598936
$ws.Cells.Item(21, 6).Value = 10134
$ws.Cells.Item(21, 7).Value = 10391
$ws.Cells.Item(21, 8).Value = -2.473294196901
$ws.Cells.Item(21, 9).Value = 84822
$ws.Cells.Item(21, 10).Value = 85026
$ws.Cells.Item(21, 11).Value = -0.239926610683
$ws.Cells.Item(21, 12).Value = 35.773854305059
$ws.Cells.Item(21, 13).Value = 22.262421263531
$ws.Cells.Item(21, 14).Value = -70.540551804617
$ws.Cells.Item(22, 3).Value = 39
$ws.Cells.Item(22, 4).Value = 44
$ws.Cells.Item(22, 5).Value = -11.363636363636
$ws.Cells.Item(22, 7).Value = 176
$ws.Cells.Item(22, 8).Value = -12.5
$ws.Cells.Item(22, 9).Value = 1468
$ws.Cells.Item(22, 10).Value = 1534
$ws.Cells.Item(22, 11).Value = -4.302477183833
$ws.Cells.Item(22, 12).Value = 41.972920696324
$ws.Cells.Item(22, 13).Value = 4.113475177304
$ws.Cells.Item(23, 3).Value = 105
$ws.Cells.Item(23, 4).Value = 127
$ws.Cells.Item(23, 5).Value = -17.322834645669
$ws.Cells.Item(23, 6).Value = 500
$ws.Cells.Item(23, 7).Value = 490
$ws.Cells.Item(23, 8).Value = 2.040816326530
$ws.Cells.Item(23, 9).Value = 4226
$ws.Cells.Item(23, 10).Value = 4086
$ws.Cells.Item(23, 11).Value = 3.426333822809
$ws.Cells.Item(23, 12).Value = 17.096148517594
$ws.Cells.Item(23, 13).Value = 50.284495021337
$ws.Cells.Item(24, 3).Value = 2215
$ws.Cells.Item(24, 4).Value = 2396
$ws.Cells.Item(24, 5).Value = -7.554257095158
$ws.Cells.Item(24, 6).Value = 9025
$ws.Cells.Item(24, 7).Value = 9774
$ws.Cells.Item(24, 8).Value = -7.663188049928
$ws.Cells.Item(24, 9).Value = 74542
$ws.Cells.Item(24, 10).Value = 77235
$ws.Cells.Item(24, 11).Value = -3.486761183401
$ws.Cells.Item(24, 12).Value = 38.107236817727
$ws.Cells.Item(24, 13).Value = 37.290726586241
$ws.Cells.Item(25, 3).Value = 915
$ws.Cells.Item(25, 4).Value = 796
$ws.Cells.Item(25, 5).Value = 14.949748743718
$ws.Cells.Item(25, 6).Value = 3495
$ws.Cells.Item(25, 7).Value = 3127
$ws.Cells.Item(25, 8).Value = 11.768468180364
$ws.Cells.Item(25, 9).Value = 29619
$ws.Cells.Item(25, 10).Value = 28127
$ws.Cells.Item(25, 11).Value = 5.304511679169
$ws.Cells.Item(25, 12).Value = 27.404507914659
$ws.Cells.Item(25, 13).Value = -5.696001018848
$ws.Cells.Item(26, 3).Value = 44
$ws.Cells.Item(26, 4).Value = 51
$ws.Cells.Item(26, 5).Value = -13.725490196078
$ws.Cells.Item(26, 6).Value = 180
$ws.Cells.Item(26, 7).Value = 196
$ws.Cells.Item(26, 8).Value = -8.163265306122
$ws.Cells.Item(26, 9).Value = 1635
$ws.Cells.Item(26, 10).Value = 1796
$ws.Cells.Item(26, 11).Value = -8.964365256124
$ws.Cells.Item(26, 12).Value = -1.208459214501
$ws.Cells.Item(27, 3).Value = 96
$ws.Cells.Item(27, 4).Value = 91
$ws.Cells.Item(27, 5).Value = 5.494505494505
$ws.Cells.Item(27, 6).Value = 415
$ws.Cells.Item(27, 7).Value = 402
$ws.Cells.Item(27, 8).Value = 3.233830845771
$ws.Cells.Item(27, 9).Value = 3559
$ws.Cells.Item(27, 10).Value = 3467
$ws.Cells.Item(27, 11).Value = 2.653591000865
$ws.Cells.Item(27, 12).Value = 10.665422885572
$ws.Cells.Item(28, 3).Value = 18
$ws.Cells.Item(28, 4).Value = 28
$ws.Cells.Item(28, 5).Value = -35.714285714285
$ws.Cells.Item(28, 6).Value = 81
$ws.Cells.Item(28, 7).Value = 125
$ws.Cells.Item(28, 8).Value = -35.2
$ws.Cells.Item(28, 9).Value = 815
$ws.Cells.Item(28, 10).Value = 1135
$ws.Cells.Item(28, 11).Value = -28.193832599118
$ws.Cells.Item(28, 12).Value = -36.526479750778
$ws.Cells.Item(28, 13).Value = -34.432823813354
$ws.Cells.Item(28, 14).Value = -79.985265225933
$ws.Cells.Item(29, 3).Value = 16
$ws.Cells.Item(29, 4).Value = 22
$ws.Cells.Item(29, 5).Value = -27.272727272727
$ws.Cells.Item(29, 6).Value = 75
$ws.Cells.Item(29, 7).Value = 104
$ws.Cells.Item(29, 8).Value = -27.884615384615
$ws.Cells.Item(29, 9).Value = 693
$ws.Cells.Item(29, 10).Value = 938
$ws.Cells.Item(29, 11).Value = -26.119402985074
$ws.Cells.Item(29, 12).Value = -35.414725069897
$ws.Cells.Item(29, 13).Value = -32.521908471275
$ws.Cells.Item(29, 14).Value = -81.122309997275
$ws.Cells.Item(30, 4).Value = 7
$ws.Cells.Item(30, 5).Value = -100
$ws.Cells.Item(30, 6).Value = 15
$ws.Cells.Item(30, 7).Value = 42
$ws.Cells.Item(30, 8).Value = -64.285714285714
$ws.Cells.Item(30, 10).Value = 458
$ws.Cells.Item(30, 11).Value = -32.532751091703
$ws.Cells.Item(30, 12).Value = -15.803814713896

# --- Row 30 "Hate Crimes" Week-to-Date 2023 column became a text "0" -----
# Build the literal text "0" in a scratch cell (format as Text so Excel
# doesn't coerce it back to a number), then bring only the VALUE into C30;
# finally restore C30's number format/style from the sibling header cell A30
# (already General/text, right aligned) so it matches style used for similar
# "N/A"-like text cells in this column.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("C30").PasteSpecial(-4163)
$ws.Range("A30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$scratch.Clear()
